$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("difficulty")

# --- Add new data rows (5, 6, 7) ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "difficulty_local_CN_3"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 14

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "difficulty_local_CN_4"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 16

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "difficulty_local_CN_5"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 18

# --- Column widths (B, C, D) ---
# (target stored widths are 32.0909.., 21.1818.., 27.3636.. ~ values below are the
# closest achievable inputs given this host's column-width pixel/MDW rounding)
$ws.Columns.Item(2).ColumnWidth = 31.428571428571427
$ws.Columns.Item(3).ColumnWidth = 20.428571428571427
$ws.Columns.Item(4).ColumnWidth = 26.714285714285715

# --- Outline grouping for the new rows ---
$ws.Rows.Item(5).OutlineLevel = 6
$ws.Rows.Item(6).OutlineLevel = 6
$ws.Rows.Item(7).OutlineLevel = 6
$ws.Columns.Item(3).OutlineLevel = 3

# --- Active selection ---
$ws.Range("F10").Select()
